$d = $word.ActiveDocument

# 1) Replace the "Use SONAR to find where to go via ANT" bullet with two new bullets:
#    "Find the open grate" and "Show-off sonar image effect inside vent"
$sonarIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext.StartsWith("Use SONAR to find where to go via ANT")) {
        $sonarIndex = $i
        break
    }
}
$sonarPara = $d.Paragraphs.Item($sonarIndex)
$sonarPara.Range.InsertParagraphAfter()
$d.Paragraphs.Item($sonarIndex).Range.Text = "Find the open grate"
$d.Paragraphs.Item($sonarIndex + 1).Range.Text = "Show-off sonar image effect inside vent"

# 2) Move the _GoBack bookmark from the end of the "Place power cell..." bullet
#    to the middle of "AR tooltip on vent indicates" (splitting it into
#    "AR tooltip on vent indic" | bookmark | "ates it leads to ...")
$full = $d.Content.Text
$idx = $full.IndexOf("AR tooltip on vent indicates")
$splitPos = $idx + ("AR tooltip on vent indic").Length
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $splitRange)
